# Add new countries to the population dropdown list (Sheet1).
# Strategy:
#   1. Append the 9 new country rows after the existing 81 data rows
#      (rows 83-91), copying the format of the last existing data row so
#      the new rows pick up the same style (s="4") as the rest of the
#      table.
#   2. Re-sort the whole data range (A2:D91) ascending by column A, which
#      is exactly what the workbook author did in Excel (Data > Sort) -
#      this both reorders the rows alphabetically and carries the
#      per-row formatting (incl. the yellow-highlighted Czech
#      Republic / Nigeria rows) along with the data.
#   3. Tidy up the view: drop the lingering A2 selection marker and make
#      sure the sheet's used-range/dimension reflects the extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the new countries -------------------------------------------------
# New country data, in alphabetical-by-name order (this also controls the
# order new entries land in the shared-string table, matching how Excel
# itself would have appended them while typing top-to-bottom).
$newCountries = @(
    @("Algeria", "DZA", "YES", 45),
    @("Bahrain", "BHR", "NO", 1),
    @("Egypt", "EGY", "YES", 30),
    @("Ghana", "GHA", "YES", 15),
    @("Hong Kong", "HKG", "NO", 1),
    @("Jamaica", "JAM", "YES", 5),
    @("Paraguay", "PRY", "YES", 20),
    @("United Arab Emirates", "ARE", "NO", 10),
    @("Zimbabwe", "ZWE", "YES", 15)
)

$lastRow = 82
$startRow = $lastRow + 1
$endRow = $startRow + $newCountries.Count - 1

# Copy the formatting of the last data row down onto the new block so the
# appended cells share the same style as the rest of the table.
$ws.Range("A$($lastRow):D$($lastRow)").Copy()
$ws.Range("A$($startRow):D$($endRow)").PasteSpecial(-4122) # xlPasteFormats

$r = $startRow
foreach ($country in $newCountries) {
    $ws.Cells.Item($r, 1).Value2 = $country[0]
    $ws.Cells.Item($r, 2).Value2 = $country[1]
    $ws.Cells.Item($r, 3).Value2 = $country[2]
    $ws.Cells.Item($r, 4).Value2 = $country[3]
    $r = $r + 1
}

# --- 2. Re-sort the full table body A2:D91 ascending by Country (col A) ---------
$ws.Range("A2:D$($endRow)").Sort($ws.Range("A2:A$($endRow)"))

# --- 3. Clean up the view --------------------------------------------------------
$ws.Range("A2").Select()
